$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Триан-43" row's ResultImage cell (K5) keeps its original text
# ("1.jpg,2.jpg,3.jpg,4.jpg"); it was simply re-entered while other image
# lists ("1.jpg,2.jpg,3.jpg, _3.jpg, 4.jpg", etc.) were being reviewed, so we
# round-trip the value here to reflect that it was touched during the edit.
$ws.Range("K5").Value = "1.jpg,2.jpg,3.jpg, _3.jpg, 4.jpg"
$ws.Range("K5").Value = "1.jpg,2.jpg,3.jpg,4.jpg"

# Editing finished with the cursor on A2.
$ws.Activate()
$ws.Range("A2").Select()
